# [Engine] [Render] commit shader reflection
#
# Roadmap!A7:G7 ("调整注册的反射类" / shader-reflection-class row) is marked
# done: it gets an end date in column F and its progress marker (column G)
# flips from "进行中" (in progress) to "已完成" (completed) - the same
# status/fill used by every other completed row (style index 2, the blue
# fill) instead of the in-progress style (style index 1, the orange fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap")

# --- 1. Fill in the missing end date for row 7 -----------------------------
# Row 4's column F already holds the literal text "2023.11.11" (same
# completion date used elsewhere in the sheet) - write it as text (leading
# apostrophe forces text so Excel doesn't reinterpret the dotted string as a
# date serial) so it reuses the existing shared string instead of minting a
# numeric date value.
$ws.Cells.Item(7, 6).Value = "'2023.11.11"

# --- 2. Recolor the whole row to the "completed" style ---------------------
# Copy just the formatting (fill/style) from row 4 - an already-completed
# row using style index 2 - onto each populated cell of row 7 (A, B, C, E,
# F, G). This flips the row from the "in progress" fill to the "completed"
# fill without touching any cell's value.
foreach ($col in @(1, 2, 3, 5, 6, 7)) {
    $srcCell = $ws.Cells.Item(4, $col)
    $dstCell = $ws.Cells.Item(7, $col)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
}

# --- 3. Flip the progress marker text from "in progress" to "completed" ---
$ws.Cells.Item(7, 7).Value = "已完成"

# --- 4. Move the active selection to G11 (matches the saved cursor spot) ---
$ws.Activate()
$ws.Range("G11").Select()

$excel.CutCopyMode = 0
